# Add shading type ("type_shade", column L) values to the Singapore
# construction archetypes on the ARCHITECTURE sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Default shading type for all building-use rows (2-19).
$defaultShade = "T1"

# Rows that use a different shading type ("T0" - no shading).
$overrides = @{
    13 = "T0"   # SWIMMING
    15 = "T0"   # PARKING
}

for ($row = 2; $row -le 19; $row++) {
    if ($overrides.ContainsKey($row)) {
        $shade = $overrides[$row]
    } else {
        $shade = $defaultShade
    }
    $ws.Cells.Item($row, 12).Value = $shade
}

$ws.Range("L1").Select()
